$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.250.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.137.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.134.15"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.670.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.134"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.98%  "
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "58.278.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.135.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "343.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.90%  "
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0935"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.57%  "
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("E39").Value = "  -4.54%  "
$ws.Range("E40").Value = "  +10.87%  "
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("E42").Value = "  +5.11%  "
$ws.Range("E43").Value = "  +2.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.174.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0264"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.265.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("E49").Value = "  +4.55%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.65%  "
